$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'62.370.16"
$ws.Cells.Item(2, 5).Value = "  +3.46%  "
$ws.Cells.Item(3, 4).Value = "'2.408.94"
$ws.Cells.Item(3, 5).Value = "  +0.86%  "
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  +0.33%  "
$ws.Cells.Item(5, 4).Value = "'572.23"
$ws.Cells.Item(5, 5).Value = "  +1.60%  "
$ws.Cells.Item(6, 4).Value = "'145.24"
$ws.Cells.Item(6, 5).Value = "  +4.99%  "
$ws.Cells.Item(7, 4).Value = "'0.998"
$ws.Cells.Item(7, 5).Value = "  -0.26%  "
$ws.Cells.Item(8, 4).Value = "'0.539"
$ws.Cells.Item(8, 5).Value = "  +0.41%  "
$ws.Cells.Item(9, 4).Value = "'2.433.57"
$ws.Cells.Item(9, 5).Value = "  +1.78%  "
$ws.Cells.Item(10, 4).Value = "'0.111"
$ws.Cells.Item(10, 5).Value = "  +5.10%  "
$ws.Cells.Item(11, 5).Value = "  +0.84%  "
$ws.Cells.Item(12, 5).Value = "  +3.74%  "
$ws.Cells.Item(13, 4).Value = "'0.351"
$ws.Cells.Item(13, 5).Value = "  +3.93%  "
$ws.Cells.Item(14, 4).Value = "'26.75"
$ws.Cells.Item(14, 5).Value = "  +4.33%  "
$ws.Cells.Item(15, 4).Value = "'0.0000179"
$ws.Cells.Item(15, 5).Value = "  +8.14%  "
$ws.Cells.Item(16, 4).Value = "'2.878.61"
$ws.Cells.Item(16, 5).Value = "  +2.18%  "
$ws.Cells.Item(17, 4).Value = "'62.182.97"
$ws.Cells.Item(17, 5).Value = "  +3.38%  "
$ws.Cells.Item(18, 4).Value = "'2.431.36"
$ws.Cells.Item(18, 5).Value = "  +1.74%  "
$ws.Cells.Item(19, 4).Value = "'7.92"
$ws.Cells.Item(19, 5).Value = "  -5.85%  "
$ws.Cells.Item(20, 4).Value = "'10.87"
$ws.Cells.Item(20, 5).Value = "  +2.54%  "
$ws.Cells.Item(21, 4).Value = "'325.98"
$ws.Cells.Item(21, 5).Value = "  +0.34%  "
$ws.Cells.Item(22, 5).Value = "  +2.65%  "
$ws.Cells.Item(23, 4).Value = "'2.02"
$ws.Cells.Item(23, 5).Value = "  +13.86%  "
$ws.Cells.Item(24, 4).Value = "'0.999"
$ws.Cells.Item(24, 5).Value = "  -0.17%  "
$ws.Cells.Item(25, 4).Value = "'65.52"
$ws.Cells.Item(25, 5).Value = "  +1.76%  "
$ws.Cells.Item(26, 4).Value = "'615.12"
$ws.Cells.Item(26, 5).Value = "  +10.90%  "
$ws.Cells.Item(27, 4).Value = "'8.34"
$ws.Cells.Item(27, 5).Value = "  +4.94%  "
$ws.Cells.Item(28, 4).Value = "'0.0₃0984"
$ws.Cells.Item(28, 5).Value = "  +9.56%  "
$ws.Cells.Item(29, 4).Value = "'2.531.78"
$ws.Cells.Item(30, 4).Value = "'8.09"
$ws.Cells.Item(30, 5).Value = "  +2.18%  "
$ws.Cells.Item(31, 5).Value = "  +8.93%  "
$ws.Cells.Item(32, 5).Value = "  +6.31%  "
$ws.Cells.Item(33, 5).Value = "  +1.78%  "
$ws.Cells.Item(34, 5).Value = "  +4.93%  "
$ws.Cells.Item(35, 4).Value = "'0.995"
$ws.Cells.Item(35, 5).Value = "  -0.27%  "
$ws.Cells.Item(36, 5).Value = "  +5.66%  "
$ws.Cells.Item(37, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(37, 4).Value = "'0.372"
$ws.Cells.Item(37, 5).Value = "  +1.54%  "
$ws.Cells.Item(38, 2).Value = "Monero"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(38, 4).Value = "'151.50"
$ws.Cells.Item(38, 5).Value = "  -1.13%  "
$ws.Cells.Item(39, 4).Value = "'5.38"
$ws.Cells.Item(39, 5).Value = "  +7.10%  "
$ws.Cells.Item(40, 4).Value = "'18.54"
$ws.Cells.Item(40, 5).Value = "  +1.58%  "
$ws.Cells.Item(41, 4).Value = "'2.73"
$ws.Cells.Item(41, 5).Value = "  +19.86%  "
$ws.Cells.Item(43, 4).Value = "'42.31"
$ws.Cells.Item(43, 5).Value = "  +2.90%  "
$ws.Cells.Item(44, 5).Value = "  -0.03%  "
$ws.Cells.Item(45, 4).Value = "'0.0₆0278"
$ws.Cells.Item(45, 5).Value = "  -0.48%  "
$ws.Cells.Item(46, 4).Value = "'143.63"
$ws.Cells.Item(46, 5).Value = "  +0.55%  "
$ws.Cells.Item(47, 5).Value = "  +2.46%  "
$ws.Cells.Item(48, 4).Value = "'20.18"
$ws.Cells.Item(48, 5).Value = "  +7.28%  "
$ws.Cells.Item(49, 4).Value = "'0.599"
$ws.Cells.Item(49, 5).Value = "  +1.87%  "
$ws.Cells.Item(50, 4).Value = "'0.0513"
$ws.Cells.Item(50, 5).Value = "  +3.19%  "
$ws.Cells.Item(51, 4).Value = "'0.0915"
$ws.Cells.Item(51, 5).Value = "  +2.14%  "
